$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "304.20"
Set-TextValue "E2" "1.47%"
Set-TextValue "D3" "35.63"
Set-TextValue "E3" "10.47%"
Set-TextValue "D4" "5.064"
Set-TextValue "E4" "1.11%"
Set-TextValue "D5" "0.07802"
Set-TextValue "E5" "1.10%"
Set-TextValue "D6" "2.255"
Set-TextValue "E6" "0.20%"
Set-TextValue "D7" "8.100"
Set-TextValue "E7" "2.04%"
Set-TextValue "D8" "4.049"
Set-TextValue "E8" "6.07%"
Set-TextValue "D9" "0.9284"
Set-TextValue "E9" "0.59%"
Set-TextValue "D10" "0.09492"
Set-TextValue "E10" "-4.08%"
Set-TextValue "D11" "0.1825"
Set-TextValue "E11" "3.62%"
Set-TextValue "D12" "0.08511"
Set-TextValue "E12" "0.72%"
Set-TextValue "D13" "0.03460"
Set-TextValue "E13" "4.63%"
Set-TextValue "D14" "0.09940"
Set-TextValue "E14" "0.69%"
Set-TextValue "D15" "0.001479"
Set-TextValue "E15" "0.31%"
Set-TextValue "D16" "0.005692"
Set-TextValue "E16" "0.39%"
Set-TextValue "D17" "3.480"
Set-TextValue "E17" "-1.85%"
Set-TextValue "E18" "-0.27%"
Set-TextValue "D19" "0.3407"
Set-TextValue "E19" "1.77%"
Set-TextValue "D20" "0.1323"
Set-TextValue "E20" "-0.95%"
Set-TextValue "D21" "4.568"
Set-TextValue "E21" "10.85%"
Set-TextValue "D22" "0.2236"
Set-TextValue "E22" "7.26%"
Set-TextValue "D23" "0.04678"
Set-TextValue "E23" "3.19%"
Set-TextValue "D24" "0.001239"
Set-TextValue "E24" "1.84%"
Set-TextValue "D25" "0.004537"
Set-TextValue "E25" "3.85%"
Set-TextValue "D26" "0.0001300"
Set-TextValue "E26" "0.61%"
Set-TextValue "E27" "-19.99%"
Set-TextValue "D39" "0.01778"
Set-TextValue "E39" "4.46%"
Set-TextValue "D40" "0.04718"
Set-TextValue "E40" "0.99%"
Set-TextValue "D41" "0.007956"
Set-TextValue "E41" "3.10%"
Set-TextValue "D42" "0.1417"
Set-TextValue "E42" "1.52%"
Set-TextValue "D43" "0.007972"
Set-TextValue "E43" "-18.43%"
Set-TextValue "D44" "0.002223"
Set-TextValue "E44" "7.19%"
Set-TextValue "D45" "0.009090"
Set-TextValue "E45" "-6.31%"
Set-TextValue "D46" "0.00006191"
Set-TextValue "E46" "2.05%"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "0.63%"
Set-TextValue "D48" "5.412"
Set-TextValue "E48" "112.14%"
Set-TextValue "D49" "0.002691"
Set-TextValue "E49" "35.33%"
Set-TextValue "D50" "0.00002101"
Set-TextValue "E50" "0.63%"
Set-TextValue "E51" "0.63%"
